# Insert a new weekly price-report row for "Fukumoto / Segunda" (Región de
# Coquimbo) immediately above the former row 128. Excel's InsertShiftDown
# semantics push the existing rows 128..144 down to 129..145 (and the sheet
# dimension/UsedRange grows to A1:T145 automatically), exactly matching the
# recorded diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(128).Insert()

$ws.Cells.Item(128, 1).Value = 1
$ws.Cells.Item(128, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(128, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(128, 4).Value = 45077
$ws.Cells.Item(128, 5).Value = 15
$ws.Cells.Item(128, 6).Value = "Fruta"
$ws.Cells.Item(128, 7).Value = 100102
$ws.Cells.Item(128, 8).Value = "Cítricos"
$ws.Cells.Item(128, 9).Value = 100102005
$ws.Cells.Item(128, 10).Value = "Naranja"
$ws.Cells.Item(128, 11).Value = "Fukumoto"
$ws.Cells.Item(128, 12).Value = "Segunda"
$ws.Cells.Item(128, 13).Value = 350
$ws.Cells.Item(128, 14).Value = 950
$ws.Cells.Item(128, 15).Value = 1000
$ws.Cells.Item(128, 16).Value = 971
$ws.Cells.Item(128, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(128, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(128, 19).Value = 971
$ws.Cells.Item(128, 20).Value = 1
